$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: __len__ -> __getitem__
$ws.Range("B2").Value = "__getitem__"

# A2: was empty -> now holds the index argument (5), with a numeric
# (#,##0) format but keeping "general" alignment
$ws.Range("A2").Value = 5
$ws.Range("A2").NumberFormat = "#,##0"
$ws.Range("A2").HorizontalAlignment = 1

# D2: was empty -> now holds the result (1)
$ws.Range("D2").Value = 1

# C4 / C5: the referenced cell name changes from "A2" to "A3"
$ws.Range("C4").Value = "A3"
$ws.Range("C5").Value = "A3"

# D4: was the numeric literal 18 -> now a text reference to cell "A2"
$ws.Range("D4").Value = "A2"
